# A new weekly price record was inserted as row 16 (pushing the former
# rows 16-28 down to 17-29). Insert a whole row at 16 so existing rows
# (and their date-format styling in column D) shift down intact, then
# populate the new row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(16).Insert()

$ws.Range("A16").Value = 6
$ws.Range("B16").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C16").Value = "Metropolitana"
$ws.Range("D16").Value = 44763
$ws.Range("E16").Value = 13
$ws.Range("F16").Value = 100112035
$ws.Range("G16").Value = "Bruselas (repollito)"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 350
$ws.Range("K16").Value = 17000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 17571
$ws.Range("N16").Value = "$/malla 15 kilos"
$ws.Range("O16").Value = "Provincia de Quillota"
$ws.Range("P16").Value = 1171
$ws.Range("Q16").Value = 15
$ws.Range("R16").Value = "Hortaliza"
